$d = $word.ActiveDocument

# The paragraph that holds the hidden "_GoBack" bookmark is the
# "End screen: No play again or exit buttons. [FAIL]" paragraph. It gets
# merged into the preceding ("AI Boats: ...") paragraph, and all of its run
# text is removed, leaving only the (now empty) bookmark behind at the end
# of the merged paragraph.

# Locate the bookmark and the paragraph containing it (rather than a
# hard-coded paragraph index, so this keeps working even if the bookmark's
# paragraph position shifts).
$b = $d.Bookmarks("_GoBack")
$endScreenPara = $d.Range($b.Start, $b.Start).Paragraphs(1)
$prevPara = $endScreenPara.Previous()

# Step 1: clear the text inside the bookmark ("End screen: ") but keep the
# bookmarkStart/bookmarkEnd markers themselves.
$insideBookmark = $d.Range($b.Start, $b.End)
$insideBookmark.Text = ""

# Step 2: clear the remaining run text in that paragraph, i.e. everything
# after the (now empty) bookmark up to (but not including) the paragraph
# mark ("No play again or exit buttons. [FAIL]").
$b = $d.Bookmarks("_GoBack")
$tail = $d.Range($b.End, $endScreenPara.Range.End - 1)
$tail.Text = ""

# Step 3: merge that now-empty paragraph into the previous paragraph by
# deleting the paragraph mark between them. Selection.Delete (rather than
# Range.Text = "") actually collapses the two <w:p> elements while
# preserving run formatting and the bookmark.
$paraMark = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
$paraMark.Select()
$word.Selection.Delete()
